$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy number-format/style of an existing date cell onto the new date cells
# so we reuse the existing style index instead of Excel creating a new one.
$ws.Range("A34").Copy()
$ws.Range("A35:A36").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 35
$ws.Range("A35").Value = "2/20/2020"
$ws.Range("B35").Value = "using data connectivity"

# Row 36
$ws.Range("A36").Value = "2/24/2020"
$ws.Range("B36").Value = "makes some changes in UI screens."

# Update the view so row 13 is at the top and B36 is the selected / active cell,
# matching the author's saved view state.
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("B36").Select()
